# Generate Report for Handoff
# Adds a new file entry (f26e92ad-55d0-40a4-bf24-14eee9e32772.md) as row 9
# to the Overview, zh-cn and de-de sheets/tables.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

$guidMd   = "f26e92ad-55d0-40a4-bf24-14eee9e32772.md"
$guidPath = "e2e\f26e92ad-55d0-40a4-bf24-14eee9e32772.md"
$zhXlf    = "f26e92ad-55d0-40a4-bf24-14eee9e32772.f7789f0a9af2e97e87c3ca8c3716e6122b1376e3.zh-cn.xlf"
$deXlf    = "f26e92ad-55d0-40a4-bf24-14eee9e32772.f7789f0a9af2e97e87c3ca8c3716e6122b1376e3.de-de.xlf"
$hoDate   = "2016-08-12 06:54:21"
$zhDate   = "2016-08-12 06:54:14"
$deDate   = "2016-08-12 06:54:21"
$dtFmt    = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet 1: Overview  (table "Overview", columns A:G)
# ---------------------------------------------------------------------------
$lo1 = $ws1.ListObjects.Item(1)
$lo1.ListRows.Add() | Out-Null

$ws1.Range("A9").Value = $guidMd
$ws1.Range("B9").Value = $guidPath
$ws1.Range("C9").Value = ".md"
$ws1.Range("D9").Value = "'"
$ws1.Range("D9").Style = "Normal"
$ws1.Range("E9").Value = "Ready for handoff"
$ws1.Range("F9").Value = "Ready for handoff"
$ws1.Range("G9").NumberFormat = $dtFmt
$ws1.Range("G9").Value = $hoDate

$ws1.Hyperlinks.Add($ws1.Range("B9"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/e2e49c1a5f2e7c5f40d86b4e5dcd5c9d4d0d7a9e/e2e/f26e92ad-55d0-40a4-bf24-14eee9e32772.md", "", "", $guidPath) | Out-Null

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn  (table "zh-cn", columns A:P)
# ---------------------------------------------------------------------------
$lo2 = $ws2.ListObjects.Item(1)
$lo2.ListRows.Add() | Out-Null

$ws2.Range("A9").Value = $guidMd
$ws2.Range("B9").Value = ".md"
$ws2.Range("C9").Value = "Ready for handoff"
$ws2.Range("D9").Value = "e2e"
$ws2.Range("E9").Value = "ht"
$ws2.Range("F9").Value = "'False"
$ws2.Range("F9").Style = "Normal"
$ws2.Range("G9").Value = $zhXlf
$ws2.Range("H9").NumberFormat = $dtFmt
$ws2.Range("H9").Value = $zhDate
$ws2.Range("I9").Value = "'"
$ws2.Range("I9").Style = "Normal"
$ws2.Range("J9").Value = "'"
$ws2.Range("J9").Style = "Normal"
$ws2.Range("K9").NumberFormat = $dtFmt
$ws2.Range("K9").Value = "0001-01-01 00:00:00"
$ws2.Range("L9").Value = "'"
$ws2.Range("L9").Style = "Normal"
$ws2.Range("M9").Value = "'True"
$ws2.Range("M9").Style = "Normal"
$ws2.Range("N9").Value = "'"
$ws2.Range("N9").Style = "Normal"
$ws2.Range("O9").Value = "'False"
$ws2.Range("O9").Style = "Normal"
$ws2.Range("P9").Value = "'"
$ws2.Range("P9").Style = "Normal"

$ws2.Hyperlinks.Add($ws2.Range("A9"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/e2e49c1a5f2e7c5f40d86b4e5dcd5c9d4d0d7a9e/e2e/f26e92ad-55d0-40a4-bf24-14eee9e32772.md", "", "", $guidMd) | Out-Null

# ---------------------------------------------------------------------------
# Sheet 3: de-de  (table "de-de", columns A:P)
# ---------------------------------------------------------------------------
$lo3 = $ws3.ListObjects.Item(1)
$lo3.ListRows.Add() | Out-Null

$ws3.Range("A9").Value = $guidMd
$ws3.Range("B9").Value = ".md"
$ws3.Range("C9").Value = "Ready for handoff"
$ws3.Range("D9").Value = "e2e"
$ws3.Range("E9").Value = "ht"
$ws3.Range("F9").Value = "'False"
$ws3.Range("F9").Style = "Normal"
$ws3.Range("G9").Value = $deXlf
$ws3.Range("H9").NumberFormat = $dtFmt
$ws3.Range("H9").Value = $deDate
$ws3.Range("I9").Value = "'"
$ws3.Range("I9").Style = "Normal"
$ws3.Range("J9").Value = "'"
$ws3.Range("J9").Style = "Normal"
$ws3.Range("K9").NumberFormat = $dtFmt
$ws3.Range("K9").Value = "0001-01-01 00:00:00"
$ws3.Range("L9").Value = "'"
$ws3.Range("L9").Style = "Normal"
$ws3.Range("M9").Value = "'True"
$ws3.Range("M9").Style = "Normal"
$ws3.Range("N9").Value = "'"
$ws3.Range("N9").Style = "Normal"
$ws3.Range("O9").Value = "'False"
$ws3.Range("O9").Style = "Normal"
$ws3.Range("P9").Value = "'"
$ws3.Range("P9").Style = "Normal"

$ws3.Hyperlinks.Add($ws3.Range("A9"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/e2e49c1a5f2e7c5f40d86b4e5dcd5c9d4d0d7a9e/e2e/f26e92ad-55d0-40a4-bf24-14eee9e32772.md", "", "", $guidMd) | Out-Null
